# Updated Instruction Softskill Assessment
#
# For every item row (2-56) on the "Booklet_FK Lagerlogistik" sheet, add the
# sub-item instruction text in column E ("Inwieweit trifft diese Aussage aud
# dich zu?" - the same wording already used in column J) and replace the old
# column J instruction with the new "Bitte klicke an." prompt (row 3 gets a
# single blank space instead, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subItemStem = "Inwieweit trifft diese Aussage aud dich zu?"
$newInstruction = "Bitte klicke an."

for ($row = 2; $row -le 56; $row++) {
    $ws.Cells.Item($row, 5).Value = $subItemStem

    if ($row -eq 3) {
        $ws.Cells.Item($row, 10).Value = " "
    } else {
        $ws.Cells.Item($row, 10).Value = $newInstruction
    }
}

$ws.Range("J3").Select() | Out-Null
